$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2, 0.267448079705191, 56.59964947207498, 0.1278662199799618, 0.602540373712727),
    @(3, 0.3149671924479802, 41.96729235941639, 2.450578998413719, 0.8130666895022006),
    @(4, 0.4717942581390715, 33.38091580073113, 3.559937583522633, 1.023593005291674),
    @(5, 0.7564205904828479, 27.75244887172457, 3.929490631381696, 1.234119321081148),
    @(6, 1.269833776201699, 23.83950404212584, 3.565748507366222, 1.444645636870622),
    @(7, 2.088060297240443, 21.00339382255697, 2.424632906425339, 1.655171952660095),
    @(8, 3.077717856651489, 18.87752502011085, 0.8290374915218717, 1.865698268449569),
    @(9, 3.529408070698828, 17.11612189554089, 0.2314275584648205, 2.076224584239043),
    @(10, 3.681295558283151, 15.63571503292816, 0.1771222661167388, 2.286750900028516),
    @(11, 3.796827757929666, 14.40276247053303, 0.15346578630132, 2.497277215817989),
    @(12, 3.893371298411551, 13.36132434918641, 0.1355355192930273, 2.707803531607463),
    @(13, 3.975382378184309, 12.47002278048386, 0.1213987975567884, 2.918329847396937),
    @(14, 4.04589578091252, 11.6985721139952, 0.1100872244624954, 3.128856163186411),
    @(15, 4.107201075838239, 11.02432787683328, 0.1008364647788515, 3.339382478975885),
    @(16, 4.161153794865354, 10.43003911632785, 0.09284531370782054, 3.549908794765359),
    @(17, 4.208935264265145, 9.902266304978705, 0.08600934679737154, 3.760435110554833),
    @(18, 4.251539985907844, 9.430432716551824, 0.08011889574567313, 3.970961426344306),
    @(19, 4.289756738139515, 9.006090572991207, 0.0750174006789331, 4.18148774213378),
    @(20, 4.324250590989138, 8.622417339452518, 0.07052326305925839, 4.392014057923253),
    @(21, 4.355492648673713, 8.273824312724505, 0.06663267261113104, 4.602540373712726),
)

foreach ($rowdata in $data) {
    $r = $rowdata[0]
    $ws.Cells.Item($r, 3).Value = $rowdata[1]
    $ws.Cells.Item($r, 4).Value = $rowdata[2]
    $ws.Cells.Item($r, 5).Value = $rowdata[3]
    $ws.Cells.Item($r, 6).Value = $rowdata[4]
}
